$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Update header/value cells that previously read "Battery Alarm (A)" / "Battery Standby (A)"
# to the new wording "Alarm Current(A)" / "Standby Current(A)"
$ws.Range("H8").Value = "Alarm Current(A)"
$ws.Range("H9").Value = "Alarm Current(A)"
$ws.Range("H10").Value = "Alarm Current(A)"

$ws.Range("I8").Value = "Standby Current(A)"
$ws.Range("I9").Value = "Standby Current(A)"
$ws.Range("I10").Value = "Standby Current(A)"

# Update the active selection to G8 (as reflected in the saved view state)
$ws.Range("G8").Select()
